# Add a "filter section for discharge measurements": duplicate the R
# column header and add a new "C" header next to the existing row-3
# labels (R / W / Vmax) on Sheet3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "R"
$ws.Range("E3").Value = "C"

# Return the selection cursor to the top-left cell (closest achievable
# state to the saved file's default view).
$ws.Range("A1").Select() | Out-Null
